$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new "Quantity" header in F1, copying the header style/format from E1
$ws.Range("F1").Value = "Quantity"
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122) # xlPasteFormats

# Set the quantity values for each data row
$ws.Range("F2").Value = 10
$ws.Range("F3").Value = 5
$ws.Range("F4").Value = 7

# Match column F width (16.0) the same way the other custom-width columns are set
$ws.Columns.Item(6).ColumnWidth = 15.166666666666666
